$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4210.476
$ws.Range("J17").Value = 4381
$ws.Range("L17").Value = 13143
$ws.Range("N17").Value = -13479

$ws.Range("H112").Value = 1624317.9
$ws.Range("J112").Value = 2067207.5
$ws.Range("L112").Value = 6201622.5
$ws.Range("N112").Value = -6203838.5

$ws.Range("H114").Value = 89629.336
$ws.Range("J114").Value = 89629.336
$ws.Range("L114").Value = 89629.336
$ws.Range("N114").Value = -98307.336

$ws.Range("H127").Value = 936.625
$ws.Range("I127").Value = 689.7273
$ws.Range("J127").Value = 1479.8
$ws.Range("K127").Value = 2069.1819
$ws.Range("L127").Value = 4439.4
$ws.Range("M127").Value = 2890.8181
$ws.Range("N127").Value = -14359.4

$ws.Range("H137").Value = 22751.152
$ws.Range("I137").Value = 701.17145
$ws.Range("J137").Value = 92910.17999999999
$ws.Range("K137").Value = 2103.51435
$ws.Range("L137").Value = 278730.54
$ws.Range("M137").Value = 446.4856499999996
$ws.Range("N137").Value = -283830.54

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1928.4
$ws.Range("I61").Value = 1509.5454
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1509.5454
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1297.5454
$ws.Range("N61").Value = -5424

$ws.Range("H80").Value = 26037.5
$ws.Range("J80").Value = 26037.5
$ws.Range("L80").Value = 26037.5
$ws.Range("N80").Value = -28033.5

$ws.Range("H83").Value = 26037.5
$ws.Range("J83").Value = 26037.5
$ws.Range("L83").Value = 78112.5
$ws.Range("N83").Value = -88096.5

$ws.Range("H97").Value = 1483.7142
$ws.Range("I97").Value = 1769.6
$ws.Range("J97").Value = 769
$ws.Range("K97").Value = 1769.6
$ws.Range("L97").Value = 769
$ws.Range("M97").Value = -1273.6
$ws.Range("N97").Value = -1761

$ws.Range("H131").Value = 59071.668
$ws.Range("J131").Value = 59071.668
$ws.Range("L131").Value = 59071.668
$ws.Range("N131").Value = -69151.66800000001

$ws.Range("H133").Value = 39568
$ws.Range("J133").Value = 39568
$ws.Range("L133").Value = 39568
$ws.Range("N133").Value = -44628

$ws.Range("H136").Value = 1928.4
$ws.Range("I136").Value = 1509.5454
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4528.6362
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1978.6362
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 27553.5
$ws.Range("J61").Value = 27553.5
$ws.Range("L61").Value = 27553.5
$ws.Range("N61").Value = -28179.5

$ws.Range("H126").Value = 37780
$ws.Range("J126").Value = 37780
$ws.Range("L126").Value = 37780
$ws.Range("N126").Value = -47660

$ws.Range("H130").Value = 28889.5
$ws.Range("J130").Value = 28889.5
$ws.Range("L130").Value = 28889.5
$ws.Range("N130").Value = -38929.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 49999.332
$ws.Range("J20").Value = 49999.332
$ws.Range("L20").Value = 49999.332
$ws.Range("N20").Value = -50471.332

$ws.Range("H30").Value = 49999.332
$ws.Range("J30").Value = 49999.332
$ws.Range("L30").Value = 49999.332
$ws.Range("N30").Value = -50181.332

$ws.Range("H128").Value = 49999.332
$ws.Range("J128").Value = 49999.332
$ws.Range("L128").Value = 49999.332
$ws.Range("N128").Value = -59959.332

$ws.Range("H132").Value = 1589.3784
$ws.Range("I132").Value = 1325.6333
$ws.Range("J132").Value = 2719.7144
$ws.Range("K132").Value = 3976.8999
$ws.Range("L132").Value = 8159.1432
$ws.Range("M132").Value = -1446.8999
$ws.Range("N132").Value = -13219.1432

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 5766.25
$ws.Range("I3").Value = 3358.0908
$ws.Range("J3").Value = 11064.2
$ws.Range("K3").Value = 10074.2724
$ws.Range("L3").Value = 33192.60000000001
$ws.Range("M3").Value = -9962.2724
$ws.Range("N3").Value = -33416.60000000001

$ws.Range("H5").Value = 954.9643
$ws.Range("J5").Value = 1399.6154
$ws.Range("L5").Value = 4198.8462
$ws.Range("N5").Value = -4422.8462

$ws.Range("H92").Value = 447.5263
$ws.Range("I92").Value = 392
$ws.Range("J92").Value = 509.22223
$ws.Range("K92").Value = 1176
$ws.Range("L92").Value = 1527.66669
$ws.Range("M92").Value = 72
$ws.Range("N92").Value = -4023.66669

$ws.Range("H131").Value = 2223117.8
$ws.Range("J131").Value = 1002.2917
$ws.Range("L131").Value = 3006.8751
$ws.Range("N131").Value = -13086.8751

$ws.Range("H132").Value = 1171.1923
$ws.Range("I132").Value = 550
$ws.Range("J132").Value = 1447.2778
$ws.Range("K132").Value = 4950
$ws.Range("L132").Value = 13025.5002
$ws.Range("M132").Value = -2420
$ws.Range("N132").Value = -18085.5002

$ws.Range("H135").Value = 954.9643
$ws.Range("J135").Value = 1399.6154
$ws.Range("L135").Value = 12596.5386
$ws.Range("N135").Value = -17666.5386

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 5250000
$ws.Range("I7").Value = 5250000
$ws.Range("K7").Value = 5250000
$ws.Range("M7").Value = -5249888

$ws.Range("H8").Value = 5250000
$ws.Range("I8").Value = 5250000
$ws.Range("K8").Value = 5250000
$ws.Range("M8").Value = -5249861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1335.2142
$ws.Range("I46").Value = 1853.5714
$ws.Range("J46").Value = 816.8570999999999
$ws.Range("K46").Value = 1853.5714
$ws.Range("L46").Value = 816.8570999999999
$ws.Range("M46").Value = -1665.5714
$ws.Range("N46").Value = -1192.8571

$ws.Range("H93").Value = 1219.75
$ws.Range("I93").Value = 1228.2778
$ws.Range("J93").Value = 1204.4
$ws.Range("K93").Value = 1228.2778
$ws.Range("L93").Value = 1204.4
$ws.Range("M93").Value = 19.72219999999993
$ws.Range("N93").Value = -3700.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 61131
$ws.Range("J46").Value = 61131
$ws.Range("L46").Value = 61131
$ws.Range("N46").Value = -61593

$ws.Range("H126").Value = 5885.6
$ws.Range("I126").Value = 7754
$ws.Range("J126").Value = 1081.1428
$ws.Range("K126").Value = 23262
$ws.Range("L126").Value = 3243.4284
$ws.Range("M126").Value = -20792
$ws.Range("N126").Value = -8183.428400000001

$ws.Range("H134").Value = 61131
$ws.Range("J134").Value = 61131
$ws.Range("L134").Value = 183393
$ws.Range("N134").Value = -188463
